# Rename the third column of the "Tabela1" table (header "Siglas") to "UF"
# on the DimEstados sheet, then leave the active cell/selection at E7 -
# matching the author's final state in the workbook after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DimEstados")

$table = $ws.ListObjects.Item("Tabela1")
$ufColumn = $table.ListColumns.Item("Siglas")
$ufColumn.Range.Item(1).Value = "UF"

$ws.Activate()
$ws.Range("E7").Select()
